# Create main workflow to run all process steps
# Update MockDocNumber values (column D) for rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D10").NumberFormat = "@"

$ws.Range("D2").Value = "1900054723"
$ws.Range("D3").Value = "1900030872"
$ws.Range("D4").Value = "1900016668"
$ws.Range("D5").Value = "1900036977"
$ws.Range("D6").Value = "1900091732"
$ws.Range("D7").Value = "1900089981"
$ws.Range("D8").Value = "1900054759"
$ws.Range("D9").Value = "1900094858"
$ws.Range("D10").Value = "1900089434"
